# ------------------------------------------------------------------
# Adds a new "2022-Q3" quarter sheet (fund-holdings detail) right
# before the existing "2022-Q2" sheet, and updates the "总计"
# (summary) sheet with a new top row for 2022-Q3, shifting the old
# rows down by one.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- 1. Update the "总计" summary sheet -------------------------
$summary = $wb.Worksheets.Item("总计")

# Existing rows 2-8 (2022-Q2 .. 2020-Q4) each slide down by one row
# to rows 3-9; a brand new row 2 is written for 2022-Q3.
$summaryRows = @(
    ,@('2022-Q3', 33, '10.23')
    ,@('2022-Q2', 25, '10.24')
    ,@('2022-Q1', 24, '9.210000000000001')
    ,@('2021-Q4', 19, '8.109999999999999')
    ,@('2021-Q3', 23, '12.57')
    ,@('2021-Q2', 26, '6.4')
    ,@('2021-Q1', 22, '6.31')
    ,@('2020-Q4', 8, '3.81')
)

# Row 9 is brand new territory for this sheet (previously A1:D8) -
# give its column-A cell the same boxed/bold/centered style used by
# the rest of column A before writing into it.
$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $summary.Cells.Item($r, 1).Value = $i
    # quarter labels ("2022-Q3", …) are never number-like, so Excel
    # keeps them as plain text without needing a quote-prefix
    $summary.Cells.Item($r, 2).Value = $row[0]
    $summary.Cells.Item($r, 3).Value = $row[1]
    $summary.Cells.Item($r, 4).Value = [double]$row[2]
}

# ---- 2. Insert the new "2022-Q3" detail sheet --------------------
$target = $wb.Worksheets.Item("2022-Q2")
$target.Copy($target)
$new = $wb.Worksheets.Item("2022-Q2 (2)")
$new.Name = "2022-Q3"

$fundData = @(
    ,@('260112', '景顺长城能源基建混合A', '22.22', '61.72', '7.88', '1.7509', 3)
    ,@('000979', '景顺长城沪港深精选股票', '20.32', '80.04', '7.90', '1.6053', 6)
    ,@('008850', '景顺长城价值稳进三年定期开放灵活配置混合', '17.29', '84.84', '8.49', '1.4679', 3)
    ,@('008715', '景顺长城价值驱动一年持有期灵活配置混合', '8.44', '90.91', '9.77', '0.8246', 4)
    ,@('010264', '鹏华成长智选混合A', '35.33', '65.51', '2.00', '0.7066', 5)
    ,@('009098', '景顺长城价值领航两年持有期混合', '7.16', '90.92', '9.32', '0.6673', 4)
    ,@('260116', '景顺长城核心竞争力混合A', '16.26', '72.95', '3.61', '0.5870', 5)
    ,@('159611', '广发中证全指电力ETF', '16.62', '99.58', '3.39', '0.5634', 8)
    ,@('008060', '景顺长城价值边际灵活配置混合A', '5.45', '80.34', '7.76', '0.4229', 7)
    ,@('009190', '景顺长城核心优选一年持有期混合', '11.52', '86.01', '3.13', '0.3606', 7)
    ,@('161123', '易方达并购重组指数（LOF）', '4.34', '94.11', '4.00', '0.1736', 5)
    ,@('673110', '西部利得新润灵活配置混合A', '4.48', '76.88', '2.84', '0.1272', 5)
    ,@('159625', '嘉实国证绿色电力ETF', '3.41', '98.77', '3.33', '0.1136', 8)
    ,@('004823', '上投摩根安裕回报混合A', '5.35', '25.71', '1.87', '0.1000', 4)
    ,@('004824', '上投摩根安裕回报混合C', '4.91', '25.71', '1.87', '0.0918', 4)
    ,@('007146', '鹏华研究智选混合', '4.05', '76.25', '1.94', '0.0786', 4)
    ,@('009490', '泰康科技创新一年定期开放混合', '2.44', '79.62', '3.16', '0.0771', 9)
    ,@('010265', '鹏华成长智选混合C', '3.12', '65.51', '2.00', '0.0624', 5)
    ,@('015779', '景顺长城价值边际灵活配置混合C', '0.79', '80.34', '7.76', '0.0613', 7)
    ,@('512390', '平安MSCI中国A股低波动ETF', '2.70', '97.87', '2.09', '0.0564', 4)
    ,@('080005', '长盛量化红利混合', '1.89', '61.68', '2.45', '0.0463', 6)
    ,@('002935', '泰康恒泰回报灵活配置混合C', '2.53', '22.01', '1.82', '0.0460', 3)
    ,@('561700', '博时中证全指电力公用事业ETF', '1.31', '98.79', '3.37', '0.0441', 8)
    ,@('005732', '富国臻选成长灵活配置混合', '2.02', '52.36', '2.10', '0.0424', 9)
    ,@('006700', '红土创新稳健混合A', '0.66', '27.52', '5.21', '0.0344', 2)
    ,@('561560', '华泰柏瑞中证全指电力公用事业ETF', '0.90', '98.29', '3.38', '0.0304', 8)
    ,@('006701', '红土创新稳健混合C', '0.50', '27.52', '5.21', '0.0260', 2)
    ,@('002934', '泰康恒泰回报灵活配置混合A', '1.12', '22.01', '1.82', '0.0204', 3)
    ,@('562350', '银华中证全指电力公用事业ETF', '0.46', '97.99', '3.34', '0.0154', 8)
    ,@('960008', '景顺长城核心竞争力混合H', '0.34', '72.95', '3.61', '0.0123', 5)
    ,@('009188', '鹏华股息精选混合', '0.62', '86.55', '1.97', '0.0122', 2)
    ,@('015731', '景顺长城核心竞争力混合C', '0.04', '72.95', '3.61', '0.0014', 5)
    ,@('015356', '西部利得新润灵活配置混合C', '0.00', '76.88', '2.84', 0, 5)
)

# Rows 27-34 do not exist yet on the copied sheet (the source quarter
# only had 25 funds) - stamp column A with the same boxed style used
# by the existing A-column cells before filling values in.
$new.Range("A2").Copy()
$new.Range("A27:A34").PasteSpecial(-4122)

for ($i = 0; $i -lt $fundData.Count; $i++) {
    $r = $i + 2
    $row = $fundData[$i]
    $new.Cells.Item($r, 1).Value = $i
    $new.Cells.Item($r, 2).Value = "'" + $row[0]
    $new.Cells.Item($r, 3).Value = $row[1]
    $new.Cells.Item($r, 4).Value = "'" + $row[2]
    $new.Cells.Item($r, 5).Value = "'" + $row[3]
    $new.Cells.Item($r, 6).Value = "'" + $row[4]
    if ($row[5] -is [string]) {
        $new.Cells.Item($r, 7).Value = "'" + $row[5]
    } else {
        $new.Cells.Item($r, 7).Value = $row[5]
    }
    $new.Cells.Item($r, 8).Value = $row[6]
}

